# ---------------------------------------------------------------------------
# "added ratio analysis for binary search"
#
# Adds a new worksheet "Binary Search" (placed before the existing
# "Merge Sort" sheet) containing a ratio-analysis table, mirroring the
# layout already used on the "Merge Sort" sheet but without the extra
# "Execution Time" input column.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$mergeSort = $wb.Worksheets("Merge Sort")

# New sheet, inserted before "Merge Sort" so tab order becomes
# [Binary Search, Merge Sort].
$ws = $wb.Worksheets.Add($mergeSort)
$ws.Name = "Binary Search"

# ---------------------------------------------------------------------------
# Column widths (character units -> Excel pads by 5/6 char when it writes
# the OOXML <col width>, so subtract that back out to land on the target).
# ---------------------------------------------------------------------------
$pad = 5 / 6
$ws.Columns.Item(2).ColumnWidth  = 14 - $pad
$ws.Columns.Item(3).ColumnWidth  = 23.6640625 - $pad
$ws.Columns.Item(4).ColumnWidth  = 9.44140625 - $pad
$ws.Columns.Item(5).ColumnWidth  = 11.5546875 - $pad
$ws.Columns.Item(6).ColumnWidth  = 11.109375 - $pad
$ws.Columns.Item(7).ColumnWidth  = 9.88671875 - $pad
$ws.Columns.Item(8).ColumnWidth  = 10.88671875 - $pad
$ws.Columns.Item(9).ColumnWidth  = 8.88671875 - $pad
$ws.Columns.Item(10).ColumnWidth = 10.5546875 - $pad
$ws.Columns.Item(11).ColumnWidth = 14 - $pad
$ws.Columns.Item(12).ColumnWidth = 10 - $pad
$ws.Columns.Item(13).ColumnWidth = 10.88671875 - $pad

# ---------------------------------------------------------------------------
# Row 2 - merged title banner
# ---------------------------------------------------------------------------
$title = $ws.Range("B2:M2")
$title.Merge()
$ws.Range("B2").Value = "Ratio Analysis of Binary Search"
$title.Font.Size = 16
$title.HorizontalAlignment = -4108   # xlCenter
$title.Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# Row 3 - column headers
# ---------------------------------------------------------------------------
$headerRow = $ws.Range("B3:M3")
$headerRow.Font.Size = 16
$headerRow.HorizontalAlignment = -4108  # xlCenter
$headerRow.VerticalAlignment = -4108    # xlCenter
$headerRow.Borders.LineStyle = 1
$ws.Rows.Item(3).RowHeight = 49.95

$ws.Range("B3").Value = "Degree(n)"
$ws.Range("C3").Value = "Number of comparisons (f(n))"
$ws.Range("D3").Value = "n"
$ws.Range("E3").Value = "log n"
$ws.Range("F3").Value = "n * log n"
$ws.Range("G3").Value = "n^2"
$ws.Range("H3").Value = "n^3"
$ws.Range("I3").Value = "f(n)/n"
$ws.Range("J3").Value = "f(n)/log n"
$ws.Range("K3").Value = "f(n)/n * log n"
$ws.Range("L3").Value = "f(n)/n^2"
$ws.Range("M3").Value = "f(n)/n^3"

# C3/K3 carry the long wrapped labels in the source sheet too - match it.
$ws.Range("C3,K3").WrapText = $true

# Superscript the trailing exponent on the n^2 / n^3 / f(n)/n^2 / f(n)/n^3
# headers, copying the rich-text run formatting already used on "Merge Sort".
$mergeSort.Range("H3").Copy() | Out-Null
$ws.Range("G3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats (n^2 header)
$mergeSort.Range("I3").Copy() | Out-Null
$ws.Range("H3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats (n^3 header)
$mergeSort.Range("J3").Copy() | Out-Null
$ws.Range("L3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats (f(n)/n^2 header)
$mergeSort.Range("K3").Copy() | Out-Null
$ws.Range("M3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats (f(n)/n^3 header)
$ws.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------------
# Rows 4-11 - data + formulas
# ---------------------------------------------------------------------------
$rows = @(
    @{ Row = 4;  B = 1;     C = 2 },
    @{ Row = 5;  B = 10;    C = 8 },
    @{ Row = 6;  B = 50;    C = 11 },
    @{ Row = 7;  B = 100;   C = 17 },
    @{ Row = 8;  B = 500;   C = 20 },
    @{ Row = 9;  B = 1000;  C = 29 },
    @{ Row = 10; B = 5000;  C = 29 },
    @{ Row = 11; B = 10000; C = 35 }
)

foreach ($r in $rows) {
    $i = $r.Row
    $ws.Range("B$i").Value = $r.B
    $ws.Range("C$i").Value = $r.C
    $ws.Range("D$i").Formula = "=B$i"
    $ws.Range("E$i").Formula = "=LOG(D$i,2)"
    $ws.Range("F$i").Formula = "=D$i*E$i"
    $ws.Range("G$i").Formula = "=POWER(D$i,2)"
    $ws.Range("H$i").Formula = "=POWER(D$i, 3)"
    $ws.Range("I$i").Formula = "=C$i/D$i"
    $ws.Range("J$i").Formula = "=C$i/E$i"
    $ws.Range("K$i").Formula = "=C$i/F$i"
    $ws.Range("L$i").Formula = "=C$i/G$i"
    $ws.Range("M$i").Formula = "=C$i/H$i"
}

$dataRange = $ws.Range("B4:M11")
$dataRange.Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# View state to mirror the source workbook
# ---------------------------------------------------------------------------
$ws.Range("C3").Select() | Out-Null
$excel.ActiveWindow.Zoom = 123

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
